$wb = $excel.ActiveWorkbook

# Update the Date value on the "Metadata" sheet
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2023-08-30T14:09:22+00:00"

# Add the Definition for the "phase-3-phase-4" concept on the "Concepts" sheet
$conceptsSheet = $wb.Worksheets.Item("Concepts")
$conceptsSheet.Range("D2").Value = "Trials that are a combination of phases III and IV."
